# Updated cryptos list values (Price column D, Volume(1h) column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.681.57'
$ws.Range('E2').Value = '  -4.80%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.209.39'
$ws.Range('E3').Value = '  -5.83%  '
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.09'
$ws.Range('E5').Value = '  +2.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.627'
$ws.Range('E6').Value = '  -5.72%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '70.49'
$ws.Range('E7').Value = '  -3.29%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -8.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.41'
$ws.Range('E10').Value = '  +11.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0946'
$ws.Range('E11').Value = '  -6.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '57.93'
$ws.Range('E12').Value = '  -4.79%  '
$ws.Range('E13').Value = '  -2.88%  '
$ws.Range('E14').Value = '  -7.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.535.53'
$ws.Range('E15').Value = '  -5.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.75'
$ws.Range('E16').Value = '  -8.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.838'
$ws.Range('E17').Value = '  -6.92%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.198.20'
$ws.Range('E18').Value = '  -5.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.605.06'
$ws.Range('E19').Value = '  -4.91%  '
$ws.Range('E20').Value = '  -6.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.71'
$ws.Range('E21').Value = '  -5.16%  '
$ws.Range('E22').Value = '  -6.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.97'
$ws.Range('E23').Value = '  -7.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.04'
$ws.Range('E24').Value = '  +11.15%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  -5.68%  '
$ws.Range('E27').Value = '  -1.91%  '
$ws.Range('E28').Value = '  -2.05%  '
$ws.Range('E29').Value = '  -5.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '169.30'
$ws.Range('E30').Value = '  -3.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.40'
$ws.Range('E31').Value = '  -7.94%  '
$ws.Range('E32').Value = '  -6.39%  '
$ws.Range('E33').Value = '  -6.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0709'
$ws.Range('E34').Value = '  -3.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.07'
$ws.Range('E35').Value = '  -4.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.60'
$ws.Range('E36').Value = '  -8.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.88'
$ws.Range('E37').Value = '  +3.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.31'
$ws.Range('E38').Value = '  +19.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.26'
$ws.Range('E39').Value = '  -5.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0271'
$ws.Range('E40').Value = '  +0.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.85'
$ws.Range('E41').Value = '  -8.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '65.10'
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.92'
$ws.Range('E43').Value = '  -1.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.83'
$ws.Range('E44').Value = '  -11.36%  '
$ws.Range('E45').Value = '  -3.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0994'
$ws.Range('E46').Value = '  -5.97%  '
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.58'
$ws.Range('E48').Value = '  +7.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.23'
$ws.Range('E49').Value = '  +7.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.17'
$ws.Range('E50').Value = '  -4.21%  '
$ws.Range('E51').Value = '  +9.85%  '
